# Rename the measure from "porcentaje-superficie-regada-sobre-sau" to
# "superficie-regada-sobre-sau" everywhere it appears on the metadata sheet
# (the machine-readable code in row 2 and the iaest-measure identifier in row 3).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "superficie-regada-sobre-sau"
$ws.Range("G3").Value = "iaest-measure:superficie-regada-sobre-sau"

# Drop the stray, content-less column M (it only ever carried cell
# formatting, no data) that trailed the real M (formerly "Municipio nombre")
# column - now column L.
$ws.Columns("M").Delete()
